# Insert a new data row at row 160 (pushing the existing rows 160-209
# down to 161-210) and populate it with a new "Cebollín" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(160).Insert()

$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44736
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 100112037
$ws.Cells.Item(160, 7).Value = "Cebollín"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 2400
$ws.Cells.Item(160, 11).Value = 1400
$ws.Cells.Item(160, 12).Value = 1600
$ws.Cells.Item(160, 13).Value = 1500
$ws.Cells.Item(160, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(160, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(160, 16).Value = 250
$ws.Cells.Item(160, 17).Value = 6
$ws.Cells.Item(160, 18).Value = "Hortaliza"
